$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with corrected values ---

# Row 209
$ws.Range("B209").Value = 47518.34
$ws.Range("C209").Value = 461.53
$ws.Range("E209").Value = 47865.87
$ws.Range("J209").Value = 59947.69

# Row 236
$ws.Range("B236").Value = 137834.57
$ws.Range("D236").Value = 1580.61
$ws.Range("E236").Value = 136563.96
$ws.Range("J236").Value = 164194.48

# Row 252
$ws.Range("B252").Value = 59908.66
$ws.Range("C252").Value = 4795.139999999999
$ws.Range("E252").Value = 64703.8
$ws.Range("J252").Value = 70228.18000000001

# Row 265
$ws.Range("F265").Value = 7553.42
$ws.Range("I265").Value = 7553.42
$ws.Range("J265").Value = 173559.93

# Row 266
$ws.Range("B266").Value = 102340.64
$ws.Range("D266").Value = 475
$ws.Range("E266").Value = 101975.64
$ws.Range("F266").Value = 5114.98
$ws.Range("I266").Value = 5114.98
$ws.Range("J266").Value = 107090.62

# --- Append new rows 268-274 ---

$newRows = @(
    @{ Row=268; A=45950; B=79957.77; C=1579.62; D=7409.63; E=74127.75999999999; F=7365.44; G=0; H=233.97; I=7131.469999999999; J=81259.23 },
    @{ Row=269; A=45951; B=68331.94; C=2979.86; D=180; E=71131.8; F=4887.9; G=0; H=0; I=4887.9; J=76019.7 },
    @{ Row=270; A=45952; B=77177.16; C=899.96; D=730; E=77347.12000000001; F=12497.3; G=0; H=0; I=12497.3; J=89844.42000000001 },
    @{ Row=271; A=45953; B=87757.27; C=745.95; D=2500; E=86003.22; F=2706.75; G=0; H=0; I=2706.75; J=88709.97 },
    @{ Row=272; A=45954; B=56310.67; C=419.4; D=0; E=56730.07; F=12673.16; G=0; H=0; I=12673.16; J=69403.23 },
    @{ Row=273; A=45955; B=201.65; C=180; D=0; E=381.65; F=0; G=0; H=0; I=0; J=381.65 },
    @{ Row=274; A=45957; B=960; C=0; D=0; E=960; F=0; G=0; H=0; I=0; J=960 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Range("A267").NumberFormat
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = "TRIGO"
}
